$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Overview")

# Update cell values
$ws.Range("A2").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.md"
$ws.Range("A3").Value = "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md"

# Rebuild hyperlinks with updated display text (targets unchanged)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/21df1751-2a26-4c19-8679-12b22d725b86.md", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/582cfb1a-645f-41e2-a5ff-9db963d3d27a.md", "", "", "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md")

$ws = $wb.Worksheets.Item("zh-cn")

# Update cell values
$ws.Range("A2").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.md"
$ws.Range("F2").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.md"
$ws.Range("D2").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf"
$ws.Range("G2").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-11 20:44:30"
$ws.Range("H2").Value = "2016-03-11 20:44:46"
$ws.Range("A3").Value = "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md"
$ws.Range("F3").Value = "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md"
$ws.Range("D3").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf"
$ws.Range("G3").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-11 20:44:30"
$ws.Range("H3").Value = "2016-03-11 20:44:46"

# Rebuild hyperlinks with updated display text (targets unchanged)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/21df1751-2a26-4c19-8679-12b22d725b86.md", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/21df1751-2a26-4c19-8679-12b22d725b86.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dddf3e89bfdd2b7cdc0b38c89b17fb0e0e1d036/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/21df1751-2a26-4c19-8679-12b22d725b86.e32f2d9a09563552f6b45587a56d141ec9681cec.zh-cn.xlf", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/70f0fcbc1e41fe9fc54a993ce1dcb2ec4130b67d/e2e/21df1751-2a26-4c19-8679-12b22d725b86.md", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9cd7691674553217eb2547d90aafa20e97d4fe32/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/21df1751-2a26-4c19-8679-12b22d725b86.e32f2d9a09563552f6b45587a56d141ec9681cec.zh-cn.xlf", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/582cfb1a-645f-41e2-a5ff-9db963d3d27a.md", "", "", "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/582cfb1a-645f-41e2-a5ff-9db963d3d27a.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dddf3e89bfdd2b7cdc0b38c89b17fb0e0e1d036/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/582cfb1a-645f-41e2-a5ff-9db963d3d27a.8f19445d09d6c0fd54db8a6edc35223fdb534180.zh-cn.xlf", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/70f0fcbc1e41fe9fc54a993ce1dcb2ec4130b67d/e2e/582cfb1a-645f-41e2-a5ff-9db963d3d27a.md", "", "", "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9cd7691674553217eb2547d90aafa20e97d4fe32/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/582cfb1a-645f-41e2-a5ff-9db963d3d27a.8f19445d09d6c0fd54db8a6edc35223fdb534180.zh-cn.xlf", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf")

$ws = $wb.Worksheets.Item("de-de")

# Update cell values
$ws.Range("A2").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.md"
$ws.Range("F2").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.md"
$ws.Range("D2").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf"
$ws.Range("G2").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf"
$ws.Range("E2").Value = "2016-03-11 20:44:33"
$ws.Range("H2").Value = "2016-03-11 20:44:52"
$ws.Range("A3").Value = "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md"
$ws.Range("F3").Value = "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md"
$ws.Range("D3").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf"
$ws.Range("G3").Value = "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf"
$ws.Range("E3").Value = "2016-03-11 20:44:33"
$ws.Range("H3").Value = "2016-03-11 20:44:52"

# Rebuild hyperlinks with updated display text (targets unchanged)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/21df1751-2a26-4c19-8679-12b22d725b86.md", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/21df1751-2a26-4c19-8679-12b22d725b86.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d54a7a3ff747f5d08b9d2b7577274501bc17809/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/21df1751-2a26-4c19-8679-12b22d725b86.e32f2d9a09563552f6b45587a56d141ec9681cec.de-de.xlf", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2c59b26d0ef44ea887c1bf3e0f14c31f5a9e4402/e2e/21df1751-2a26-4c19-8679-12b22d725b86.md", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/99908dcd03ce77c9b0ab1d544e705ded9b0557bd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/21df1751-2a26-4c19-8679-12b22d725b86.e32f2d9a09563552f6b45587a56d141ec9681cec.de-de.xlf", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/582cfb1a-645f-41e2-a5ff-9db963d3d27a.md", "", "", "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/582cfb1a-645f-41e2-a5ff-9db963d3d27a.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d54a7a3ff747f5d08b9d2b7577274501bc17809/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/582cfb1a-645f-41e2-a5ff-9db963d3d27a.8f19445d09d6c0fd54db8a6edc35223fdb534180.de-de.xlf", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2c59b26d0ef44ea887c1bf3e0f14c31f5a9e4402/e2e/582cfb1a-645f-41e2-a5ff-9db963d3d27a.md", "", "", "ffffba80ea05-eea9-4c85-8ac2-b3d2bb95c076.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/99908dcd03ce77c9b0ab1d544e705ded9b0557bd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/582cfb1a-645f-41e2-a5ff-9db963d3d27a.8f19445d09d6c0fd54db8a6edc35223fdb534180.de-de.xlf", "", "", "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf")
